$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint Backlog: the two remaining (previously unassigned) backlog rows
# -- "60152 patterns" (row 16) and "60152 code smells" (row 17) -- get
# assigned to "Joana" (a new team member / new shared string).
$ws.Range("C16").Value = "Joana"
$ws.Range("C17").Value = "Joana"

# Leave the cursor where the author left it when the file was saved.
[void]$ws.Range("C17").Select()
